$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: becomes Athletic Club's updated stats (one more match played - a win)
$ws.Range("B9").Value = "Athletic Club"
$ws.Range("C9").Value = 25
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 11
$ws.Range("G9").Value = 29
$ws.Range("H9").Value = 35
$ws.Range("I9").Value = -6
$ws.Range("J9").Value = 34
$ws.Range("K9").Value = 1.36
$ws.Range("L9").Value = "L D W W W"
$ws.Range("M9").Value = 48000
$ws.Range("N9").Value = "Gorka Guruzeta - 6"
$ws.Range("O9").Value = "Unai Sim" + [char]0x00F3 + "n"

# Row 10: becomes Real Sociedad's (previous) stats, unchanged from before
$ws.Range("B10").Value = "Real Sociedad"
$ws.Range("C10").Value = 24
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 9
$ws.Range("G10").Value = 34
$ws.Range("H10").Value = 35
$ws.Range("I10").Value = -1
$ws.Range("J10").Value = 31
$ws.Range("K10").Value = 1.29
$ws.Range("L10").Value = "W W D W L"
$ws.Range("M10").Value = 31242
$ws.Range("N10").Value = "Mikel Oyarzabal - 10"
$ws.Range("O10").Value = [char]0x00C1 + "lex Remiro"

# Row 17: Elche's updated stats (one more match played - a loss)
$ws.Range("C17").Value = 25
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 37
$ws.Range("I17").Value = -5
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = "L L L D L"
$ws.Range("N17").Value = "Andr" + [char]0x00E9 + " Silva, Rafa Mir - 6"
